# "add remove liquidity via unstake"
#
# Adds a "Price" column to the CalculationAmountAndPrice sheet and fixes the
# running-balance formulas in rows 5-7 so they correctly add/remove
# liquidity for BUY/SELL rows instead of blindly re-using the add-liquidity
# pattern. Also nudges the saved selection on all three sheets and the
# sheet1 print setup, matching the authored workbook.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # CalculationAmountAndPrice
$ws2 = $wb.Worksheets.Item(2)   # CalculationLiquidity
$ws3 = $wb.Worksheets.Item(3)   # Staking

# ---------------------------------------------------------------------
# Sheet1: CalculationAmountAndPrice
# ---------------------------------------------------------------------

# New "Price" header in N1, matching the header style used by M1.
$ws1.Range("M1").Copy() | Out-Null
$ws1.Range("N1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws1.Range("N1").Value = "Price"

# New Price column values: swap amount received divided by amount traded.
$ws1.Range("N2").Formula = "=M2/E2"
$ws1.Range("N3").Formula = "=M3/E3"
$ws1.Range("N4").Formula = "=(M4/E4)/1"
$ws1.Range("N5").Formula = "=(M5/E5)/1"
$ws1.Range("N6").Formula = "=M6/E6"

# Fix the running pool-balance formulas for the BUY/SELL rows: liquidity
# removed via an "unstake" must be subtracted (not added as if it were
# newly deposited), and the counter asset must reflect what the trader
# actually paid in (M, fee H and deducted L amounts), not a copy of the
# add-liquidity math.
$ws1.Range("A5").Formula = "=A4-M4"
$ws1.Range("B5").Formula = "=B4+M4+H4"
$ws1.Range("A6").Formula = "=A5-M5"
$ws1.Range("B6").Formula = "=B5+L5+H5"
$ws1.Range("A7").Formula = "=A6+L6+H6"
$ws1.Range("B7").Formula = "=B6-M6"

# Print setup for sheet1 (A4 portrait).
$ws1.PageSetup.PaperSize = 9     # xlPaperA4
$ws1.PageSetup.Orientation = 1   # xlPortrait

# ---------------------------------------------------------------------
# Selections (cosmetic, matches the saved cursor position per sheet)
# ---------------------------------------------------------------------

$ws2.Activate()
$ws2.Range("H21").Select() | Out-Null

$ws3.Activate()
$ws3.Range("C28").Select() | Out-Null

$ws1.Activate()
$ws1.Range("C26").Select() | Out-Null
